$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The KPI product-grouping export template's placeholder bindings were
# restructured: "ProductGroupings" is now nested one level deeper, under a
# new "KpiProductGroupings" collection whose per-grouping rows live in
# "Contents" (fixing the import path).
$ws.Range("C6").Value = "{{KpiProductGroupings.Employees.KpiProductGroupings.KpiProductGroupingTypeName}}"
$ws.Range("D6").Value = "{{KpiProductGroupings.Employees.KpiProductGroupings.Contents.Code}}"
$ws.Range("E6").Value = "{{KpiProductGroupings.Employees.KpiProductGroupings.Contents.Name}}"
$ws.Range("F6").Value = "{{KpiProductGroupings.Employees.KpiProductGroupings.Contents.ItemCount}}"
$ws.Range("G6").Value = "{{KpiProductGroupings.Employees.KpiProductGroupings.Contents.Items.Code}}"
$ws.Range("H6").Value = "{{KpiProductGroupings.Employees.KpiProductGroupings.Contents.Items.Name}}"
$ws.Range("I6").Value = "{{KpiProductGroupings.Employees.KpiProductGroupings.Contents.Items.IndirectRevenue}}"
$ws.Range("J6").Value = "{{KpiProductGroupings.Employees.KpiProductGroupings.Contents.Items.IndirectStoreCounter}}"

# The new placeholders in columns G/H are considerably longer, so widen
# those two columns to fit the new text instead of sharing one fixed width.
$null = $ws.Columns.Item(7).AutoFit()
$null = $ws.Columns.Item(8).AutoFit()
$ws.Columns.Item(7).ColumnWidth = 72.6
$ws.Columns.Item(8).ColumnWidth = 73.43

# Leave the selection on the cell that was actually edited last.
$null = $ws.Range("G6").Select()
